$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns G:K, F1 text changes ---
$ws.Range("F1").Value = "instancia_input"
$ws.Range("G1").Value = "reg1"
$ws.Range("H1").Value = "Creditos"
$ws.Range("I1").Value = "periodo_idx"
$ws.Range("J1").Value = "acred1"
$ws.Range("K1").Value = "lugar"

# --- Data row (row 2): new columns G:K ---
$ws.Range("G2").Value = "014374 del 11 de diciembre de 2019"
$ws.Range("H2").Value = 160
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "Resolución 4293 del 2016"
$ws.Range("K2").Value = "Medellín"

# --- Formatting for G2: Arial font, medium black border, wrapped/centered text ---
$g2 = $ws.Range("G2")
$g2.Font.Name = "Arial"
$g2.Borders.Weight = -4138
$g2.Borders.Color = 0
$g2.HorizontalAlignment = -4131
$g2.VerticalAlignment = -4108
$g2.WrapText = $true

# --- Row heights / thick bottom border on rows 1 and 2 ---
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 56.5

# --- Sheet view: zoom / top-left cell / selection ---
$ws.Application.ActiveWindow.Zoom = 70
$ws.Application.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("I15").Select()
